$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (23-27) need the duration number format explicitly applied to
# column R, since they have no pre-existing style to inherit from.
$ws.Range("R23:R27").NumberFormat = "[hh]:mm:ss"

$ws.Range("A2").Value = 251455
$ws.Range("B2").Value = 'BIMEC 2'
$ws.Range("C2").Value = 19
$ws.Range("D2").Value = 82.765625
$ws.Range("E2").Value = '2025-05-07 07:00:00'
$ws.Range("F2").Value = '2025-05-07 07:19:00'
$ws.Range("G2").Value = '2025-05-07 07:19:00'
$ws.Range("H2").Value = '2025-05-07 08:41:45'
$ws.Range("I2").Value = 5297
$ws.Range("J2").Value = 'bobina'
$ws.Range("K2").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = 39749
$ws.Range("O2").Value = 'X'
$ws.Range("P2").Value = 39749
$ws.Range("Q2").Value = '2025-04-15 00:00:00'
$ws.Range("R2").Value = -0.3623372395833334
$ws.Range("S2").Value = 7

$ws.Range("A3").Value = 251395
$ws.Range("B3").Value = 'BIMEC 2'
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = 35.34375
$ws.Range("E3").Value = '2025-05-07 08:41:45'
$ws.Range("F3").Value = '2025-05-07 09:00:45'
$ws.Range("G3").Value = '2025-05-07 09:00:45'
$ws.Range("H3").Value = '2025-05-07 09:36:06'
$ws.Range("I3").Value = 2262
$ws.Range("J3").Value = 'bobina'
$ws.Range("K3").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 39749
$ws.Range("O3").Value = 'X'
$ws.Range("P3").Value = 39749
$ws.Range("Q3").Value = '2025-04-23 00:00:00'
$ws.Range("R3").Value = -0.4000759548611111
$ws.Range("S3").Value = 7

$ws.Range("A4").Value = 251391
$ws.Range("B4").Value = 'BIMEC 2'
$ws.Range("C4").Value = 17
$ws.Range("D4").Value = 91.640625
$ws.Range("E4").Value = '2025-05-07 09:36:06'
$ws.Range("F4").Value = '2025-05-07 09:53:06'
$ws.Range("G4").Value = '2025-05-07 09:53:06'
$ws.Range("H4").Value = '2025-05-07 11:24:45'
$ws.Range("I4").Value = 5865
$ws.Range("J4").Value = 'bobina'
$ws.Range("K4").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 70
$ws.Range("N4").Value = 39749
$ws.Range("O4").Value = 'X'
$ws.Range("P4").Value = 39749
$ws.Range("Q4").Value = '2025-04-23 00:00:00'
$ws.Range("R4").Value = -0.4755208333333333
$ws.Range("S4").Value = 7

$ws.Range("A5").Value = 251371
$ws.Range("B5").Value = 'BIMEC 2'
$ws.Range("C5").Value = 17
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = '2025-05-07 11:24:45'
$ws.Range("F5").Value = '2025-05-07 11:41:45'
$ws.Range("G5").Value = '2025-05-07 11:41:45'
$ws.Range("H5").Value = '2025-05-07 11:41:45'
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 'bobina'
$ws.Range("K5").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = '39666 (esterno)'
$ws.Range("O5").Value = 'X'
$ws.Range("P5").Value = 39666
$ws.Range("Q5").Value = '2025-04-24 00:00:00'
$ws.Range("R5").Value = -13.48732638888889
$ws.Range("S5").Value = 7

$ws.Range("A6").Value = 251396
$ws.Range("B6").Value = 'BIMEC 2'
$ws.Range("C6").Value = 19
$ws.Range("D6").Value = 35.34375
$ws.Range("E6").Value = '2025-05-07 11:41:45'
$ws.Range("F6").Value = '2025-05-07 12:00:45'
$ws.Range("G6").Value = '2025-05-07 12:00:45'
$ws.Range("H6").Value = '2025-05-07 12:36:05'
$ws.Range("I6").Value = 2262
$ws.Range("J6").Value = 'bobina'
$ws.Range("K6").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 70
$ws.Range("N6").Value = 39749
$ws.Range("O6").Value = 'X'
$ws.Range("P6").Value = 39749
$ws.Range("Q6").Value = '2025-05-02 00:00:00'
$ws.Range("R6").Value = -0.5250651041666666
$ws.Range("S6").Value = 7

$ws.Range("A7").Value = 251548
$ws.Range("B7").Value = 'BIMEC 2'
$ws.Range("C7").Value = 19
$ws.Range("D7").Value = 206.90625
$ws.Range("E7").Value = '2025-05-07 12:36:05'
$ws.Range("F7").Value = '2025-05-07 12:55:05'
$ws.Range("G7").Value = '2025-05-07 12:55:05'
$ws.Range("H7").Value = '2025-05-08 08:22:00'
$ws.Range("I7").Value = 13242
$ws.Range("J7").Value = 'bobina'
$ws.Range("K7").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 39749
$ws.Range("O7").Value = 'X'
$ws.Range("P7").Value = 39749
$ws.Range("Q7").Value = '2025-05-06 00:00:00'
$ws.Range("R7").Value = -1.348611111111111
$ws.Range("S7").Value = 7

$ws.Range("A8").Value = 250923
$ws.Range("B8").Value = 'BIMEC 2'
$ws.Range("C8").Value = 32
$ws.Range("D8").Value = 109.46875
$ws.Range("E8").Value = '2025-05-08 08:22:00'
$ws.Range("F8").Value = '2025-05-08 08:54:00'
$ws.Range("G8").Value = '2025-05-08 08:54:00'
$ws.Range("H8").Value = '2025-05-08 10:43:28'
$ws.Range("I8").Value = 7006
$ws.Range("J8").Value = 'bobina'
$ws.Range("K8").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9'
$ws.Range("L8").Value = 5
$ws.Range("M8").Value = 76
$ws.Range("N8").Value = 39749
$ws.Range("O8").Value = 'X'
$ws.Range("P8").Value = 39749
$ws.Range("Q8").Value = '2025-04-07 00:00:00'
$ws.Range("R8").Value = -1.446853298611111
$ws.Range("S8").Value = 2

$ws.Range("A9").Value = 251477
$ws.Range("B9").Value = 'BIMEC 2'
$ws.Range("C9").Value = 19
$ws.Range("D9").Value = 468.734375
$ws.Range("E9").Value = '2025-05-08 10:43:28'
$ws.Range("F9").Value = '2025-05-08 11:02:28'
$ws.Range("G9").Value = '2025-05-08 11:02:28'
$ws.Range("H9").Value = '2025-05-09 10:51:12'
$ws.Range("I9").Value = 29999
$ws.Range("J9").Value = 'bobina'
$ws.Range("K9").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9'
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 76
$ws.Range("N9").Value = 39760
$ws.Range("O9").Value = 'X'
$ws.Range("P9").Value = 39760
$ws.Range("Q9").Value = '2025-04-28 00:00:00'
$ws.Range("R9").Value = -2.452224392361111
$ws.Range("S9").Value = 1

$ws.Range("A10").Value = 251651
$ws.Range("B10").Value = 'BIMEC 4'
$ws.Range("C10").Value = 29
$ws.Range("D10").Value = 767.7049180327868
$ws.Range("E10").Value = '2025-05-09 07:00:00'
$ws.Range("F10").Value = '2025-05-09 07:29:00'
$ws.Range("G10").Value = '2025-05-09 07:29:00'
$ws.Range("H10").Value = '2025-05-12 12:16:42'
$ws.Range("I10").Value = 46830
$ws.Range("J10").Value = 'bobina'
$ws.Range("K10").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R3 ;R6 ;R9'
$ws.Range("L10").Value = 5
$ws.Range("M10").Value = 76
$ws.Range("N10").Value = 39755
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = '2025-05-12 00:00:00'
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 4

$ws.Range("A11").Value = 251547
$ws.Range("B11").Value = 'BIMEC 5'
$ws.Range("C11").Value = 34
$ws.Range("D11").Value = 184.9154929577465
$ws.Range("E11").Value = '2025-05-08 07:00:00'
$ws.Range("F11").Value = '2025-05-08 07:34:00'
$ws.Range("G11").Value = '2025-05-08 07:34:00'
$ws.Range("H11").Value = '2025-05-08 10:38:54'
$ws.Range("I11").Value = 13129
$ws.Range("J11").Value = 'bobina'
$ws.Range("K11").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L11").Value = 4
$ws.Range("M11").Value = 70
$ws.Range("N11").Value = 39749
$ws.Range("O11").Value = 'X'
$ws.Range("P11").Value = 39749
$ws.Range("Q11").Value = '2025-05-06 00:00:00'
$ws.Range("R11").Value = -1.443691314548611
$ws.Range("S11").Value = 7

$ws.Range("A12").Value = 250759
$ws.Range("B12").Value = 'BIMEC 5'
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 118.2816901408451
$ws.Range("E12").Value = '2025-05-08 10:38:54'
$ws.Range("F12").Value = '2025-05-08 11:08:54'
$ws.Range("G12").Value = '2025-05-08 11:08:54'
$ws.Range("H12").Value = '2025-05-08 13:07:11'
$ws.Range("I12").Value = 8398
$ws.Range("J12").Value = 'bobina'
$ws.Range("K12").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12'
$ws.Range("L12").Value = 4
$ws.Range("M12").Value = 76
$ws.Range("N12").Value = 39747
$ws.Range("O12").Value = 'X'
$ws.Range("P12").Value = 39747
$ws.Range("Q12").Value = '2025-03-13 00:00:00'
$ws.Range("R12").Value = -0.5466647104861111
$ws.Range("S12").Value = 1

$ws.Range("A13").Value = 250894
$ws.Range("B13").Value = 'BIMEC 5'
$ws.Range("C13").Value = 17
$ws.Range("D13").Value = 623.4084507042254
$ws.Range("E13").Value = '2025-05-08 13:07:11'
$ws.Range("F13").Value = '2025-05-08 13:24:11'
$ws.Range("G13").Value = '2025-05-08 13:24:11'
$ws.Range("H13").Value = '2025-05-12 07:47:36'
$ws.Range("I13").Value = 44262
$ws.Range("J13").Value = 'bobina'
$ws.Range("K13").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R3 ;R6 ;R9'
$ws.Range("L13").Value = 5
$ws.Range("M13").Value = 76
$ws.Range("N13").Value = 39755
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = '2025-05-05 00:00:00'
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 4

$ws.Range("A14").Value = 251284
$ws.Range("B14").Value = 'CASON'
$ws.Range("C14").Value = 40.5
$ws.Range("D14").Value = 297.0909090909091
$ws.Range("E14").Value = '2025-05-09 07:00:00'
$ws.Range("F14").Value = '2025-05-09 07:40:30'
$ws.Range("G14").Value = '2025-05-09 07:40:30'
$ws.Range("H14").Value = '2025-05-09 12:37:35'
$ws.Range("I14").Value = 16340
$ws.Range("J14").Value = 'bobina'
$ws.Range("K14").Value = 'CASON ;R6'
$ws.Range("L14").Value = 7
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = 39747
$ws.Range("O14").Value = 'X'
$ws.Range("P14").Value = 39747
$ws.Range("Q14").Value = '2025-05-12 00:00:00'
$ws.Range("R14").Value = -1.526104797974537
$ws.Range("S14").Value = 1

$ws.Range("A15").Value = 251742
$ws.Range("B15").Value = 'CASON'
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 134.8524590163935
$ws.Range("E15").Value = '2025-05-08 07:00:00'
$ws.Range("F15").Value = '2025-05-08 07:30:00'
$ws.Range("G15").Value = '2025-05-08 07:30:00'
$ws.Range("H15").Value = '2025-05-08 09:44:51'
$ws.Range("I15").Value = 8226
$ws.Range("J15").Value = 'bobina'
$ws.Range("K15").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L15").Value = 4
$ws.Range("M15").Value = 70
$ws.Range("N15").Value = 39749
$ws.Range("O15").Value = 'X'
$ws.Range("P15").Value = 39749
$ws.Range("Q15").Value = '2025-05-15 00:00:00'
$ws.Range("R15").Value = -1.406147540983796
$ws.Range("S15").Value = 7

$ws.Range("A16").Value = 251840
$ws.Range("B16").Value = 'R10'
$ws.Range("C16").Value = 25
$ws.Range("D16").Value = 93.67213114754098
$ws.Range("E16").Value = '2025-05-08 09:44:51'
$ws.Range("F16").Value = '2025-05-08 10:09:51'
$ws.Range("G16").Value = '2025-05-08 10:09:51'
$ws.Range("H16").Value = '2025-05-08 11:43:31'
$ws.Range("I16").Value = 5714
$ws.Range("J16").Value = 'bobina'
$ws.Range("K16").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L16").Value = 5
$ws.Range("M16").Value = 70
$ws.Range("N16").Value = 39758
$ws.Range("O16").Value = 'X'
$ws.Range("P16").Value = 39758
$ws.Range("Q16").Value = '2025-05-15 00:00:00'
$ws.Range("R16").Value = -0.4885587431712963
$ws.Range("S16").Value = 1

$ws.Range("A17").Value = 251229
$ws.Range("B17").Value = 'R10'
$ws.Range("C17").Value = 25
$ws.Range("D17").Value = 307.1967213114754
$ws.Range("E17").Value = '2025-05-08 11:43:31'
$ws.Range("F17").Value = '2025-05-08 12:08:31'
$ws.Range("G17").Value = '2025-05-08 12:08:31'
$ws.Range("H17").Value = '2025-05-09 09:15:43'
$ws.Range("I17").Value = 18739
$ws.Range("J17").Value = 'bobina'
$ws.Range("K17").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R9'
$ws.Range("L17").Value = 6
$ws.Range("M17").Value = 70
$ws.Range("N17").Value = '39723 (esterno)'
$ws.Range("O17").Value = 'X'
$ws.Range("P17").Value = 39723
$ws.Range("Q17").Value = '2025-05-15 00:00:00'
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 1

$ws.Range("A18").Value = 251225
$ws.Range("B18").Value = 'R12'
$ws.Range("C18").Value = 19
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = '2025-05-08 12:00:00'
$ws.Range("F18").Value = '2025-05-08 12:19:00'
$ws.Range("G18").Value = '2025-05-08 12:19:00'
$ws.Range("H18").Value = '2025-05-08 12:19:00'
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 'bobina'
$ws.Range("K18").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9'
$ws.Range("L18").Value = 4
$ws.Range("M18").Value = 76
$ws.Range("N18").Value = 39747
$ws.Range("O18").Value = 'X'
$ws.Range("P18").Value = 39747
$ws.Range("Q18").Value = '2025-04-30 00:00:00'
$ws.Range("R18").Value = -0.5131944444444444
$ws.Range("S18").Value = 1

$ws.Range("A19").Value = 251227
$ws.Range("B19").Value = 'R12'
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = '2025-05-08 12:19:00'
$ws.Range("F19").Value = '2025-05-08 12:34:00'
$ws.Range("G19").Value = '2025-05-08 12:34:00'
$ws.Range("H19").Value = '2025-05-08 12:34:00'
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 'bobina'
$ws.Range("K19").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9'
$ws.Range("L19").Value = 4
$ws.Range("M19").Value = 76
$ws.Range("N19").Value = 39746
$ws.Range("O19").Value = 'X'
$ws.Range("P19").Value = 39746
$ws.Range("Q19").Value = '2025-05-05 00:00:00'
$ws.Range("R19").Value = -2.523611111111111
$ws.Range("S19").Value = 1

$ws.Range("A20").Value = 251782
$ws.Range("B20").Value = 'R12'
$ws.Range("C20").Value = 17
$ws.Range("D20").Value = 170.0422535211268
$ws.Range("E20").Value = '2025-05-08 12:34:00'
$ws.Range("F20").Value = '2025-05-08 12:51:00'
$ws.Range("G20").Value = '2025-05-08 12:51:00'
$ws.Range("H20").Value = '2025-05-09 07:41:02'
$ws.Range("I20").Value = 12073
$ws.Range("J20").Value = 'bobina'
$ws.Range("K20").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9'
$ws.Range("L20").Value = 3
$ws.Range("M20").Value = 76
$ws.Range("N20").Value = 39754
$ws.Range("O20").Value = 'X'
$ws.Range("P20").Value = 39754
$ws.Range("Q20").Value = '2025-05-16 00:00:00'
$ws.Range("R20").Value = -0.3201682316087963
$ws.Range("S20").Value = 1

$ws.Range("A21").Value = 251164
$ws.Range("B21").Value = 'R3'
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 204.0816326530612
$ws.Range("E21").Value = '2025-05-08 07:00:00'
$ws.Range("F21").Value = '2025-05-08 07:55:00'
$ws.Range("G21").Value = '2025-05-08 07:55:00'
$ws.Range("H21").Value = '2025-05-08 11:19:04'
$ws.Range("I21").Value = 10000
$ws.Range("J21").Value = 'bobina'
$ws.Range("K21").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L21").Value = 6
$ws.Range("M21").Value = 70
$ws.Range("N21").Value = 39749
$ws.Range("O21").Value = 'X'
$ws.Range("P21").Value = 39749
$ws.Range("Q21").Value = '2025-04-22 00:00:00'
$ws.Range("R21").Value = -1.471584467118056
$ws.Range("S21").Value = 4

$ws.Range("A22").Value = 251456
$ws.Range("B22").Value = 'R3'
$ws.Range("C22").Value = 50
$ws.Range("D22").Value = 183.6530612244898
$ws.Range("E22").Value = '2025-05-08 11:19:04'
$ws.Range("F22").Value = '2025-05-08 12:09:04'
$ws.Range("G22").Value = '2025-05-08 12:09:04'
$ws.Range("H22").Value = '2025-05-09 07:12:44'
$ws.Range("I22").Value = 8999
$ws.Range("J22").Value = 'bobina'
$ws.Range("K22").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L22").Value = 3
$ws.Range("M22").Value = 70
$ws.Range("N22").Value = 39746
$ws.Range("O22").Value = 'X'
$ws.Range("P22").Value = 39746
$ws.Range("Q22").Value = '2025-05-09 00:00:00'
$ws.Range("R22").Value = -3.300510204085648
$ws.Range("S22").Value = 1

$ws.Range("A23").Value = 251416
$ws.Range("B23").Value = 'R3'
$ws.Range("C23").Value = 40
$ws.Range("D23").Value = 229.0204081632653
$ws.Range("E23").Value = '2025-05-09 07:12:44'
$ws.Range("F23").Value = '2025-05-09 07:52:44'
$ws.Range("G23").Value = '2025-05-09 07:52:44'
$ws.Range("H23").Value = '2025-05-09 11:41:45'
$ws.Range("I23").Value = 11222
$ws.Range("J23").Value = 'bobina'
$ws.Range("K23").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L23").Value = 2
$ws.Range("M23").Value = 70
$ws.Range("N23").Value = 39755
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = '2025-04-23 00:00:00'
$ws.Range("R23").Value = 0
$ws.Range("S23").Value = 2

$ws.Range("A24").Value = 251050
$ws.Range("B24").Value = 'R6'
$ws.Range("C24").Value = 217
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = '2025-05-09 07:00:00'
$ws.Range("F24").Value = '2025-05-09 10:37:00'
$ws.Range("G24").Value = '2025-05-09 10:37:00'
$ws.Range("H24").Value = '2025-05-09 10:37:00'
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 'bobina'
$ws.Range("K24").Value = 'R6'
$ws.Range("L24").Value = 38
$ws.Range("M24").Value = 70
$ws.Range("N24").Value = 39747
$ws.Range("O24").Value = 'X'
$ws.Range("P24").Value = 39747
$ws.Range("Q24").Value = '2025-04-16 00:00:00'
$ws.Range("R24").Value = -1.442361111111111
$ws.Range("S24").Value = 2

$ws.Range("A25").Value = 251054
$ws.Range("B25").Value = 'R6'
$ws.Range("C25").Value = 35
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = '2025-05-09 10:37:00'
$ws.Range("F25").Value = '2025-05-09 11:12:00'
$ws.Range("G25").Value = '2025-05-09 11:12:00'
$ws.Range("H25").Value = '2025-05-09 11:12:00'
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 'bobina'
$ws.Range("K25").Value = 'R6'
$ws.Range("L25").Value = 38
$ws.Range("M25").Value = 70
$ws.Range("N25").Value = 39747
$ws.Range("O25").Value = 'X'
$ws.Range("P25").Value = 39747
$ws.Range("Q25").Value = '2025-04-16 00:00:00'
$ws.Range("R25").Value = -1.466666666666667
$ws.Range("S25").Value = 2

$ws.Range("A26").Value = 251081
$ws.Range("B26").Value = 'R6'
$ws.Range("C26").Value = 125
$ws.Range("D26").Value = 42.42253521126761
$ws.Range("E26").Value = '2025-05-09 11:12:00'
$ws.Range("F26").Value = '2025-05-09 13:17:00'
$ws.Range("G26").Value = '2025-05-09 13:17:00'
$ws.Range("H26").Value = '2025-05-09 13:59:25'
$ws.Range("I26").Value = 3012
$ws.Range("J26").Value = 'bobina'
$ws.Range("K26").Value = 'R6'
$ws.Range("L26").Value = 20
$ws.Range("M26").Value = 70
$ws.Range("N26").Value = '39750 (esterno)'
$ws.Range("O26").Value = 'X'
$ws.Range("P26").Value = 39750
$ws.Range("Q26").Value = '2025-04-23 00:00:00'
$ws.Range("R26").Value = -16.58293231612268
$ws.Range("S26").Value = 2

$ws.Range("A27").Value = 251706
$ws.Range("B27").Value = 'T3'
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 50.79365079365079
$ws.Range("E27").Value = '2025-05-12 07:00:00'
$ws.Range("F27").Value = '2025-05-12 07:00:00'
$ws.Range("G27").Value = '2025-05-12 07:00:00'
$ws.Range("H27").Value = '2025-05-12 07:50:47'
$ws.Range("I27").Value = 3200
$ws.Range("J27").Value = 'foglio'
$ws.Range("K27").Value = 'T3'
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = '39764 (esterno)'
$ws.Range("O27").Value = 'X'
$ws.Range("P27").Value = 39764
$ws.Range("Q27").Value = '2025-05-14 00:00:00'
$ws.Range("R27").Value = 0
$ws.Range("S27").Value = 7
